$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "DC-01"
$ws.Range("A4").Value = "DC-02"
$ws.Range("D4").Value = "Modelo"
$ws.Range("A5").Value = "DC-03"
